# feat: add 2022-Q1 data
#
# 1. Insert a new worksheet "2022-Q1" right after "2021-Q4" (before "总计"),
#    cloned from "2021-Q4" so it inherits the same header/body styling.
# 2. Fill it in with the 2022-Q1 per-fund holdings (3 funds).
# 3. Insert a matching summary row at the top of "总计"'s data and
#    renumber the existing index column.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Clone "2021-Q4" (3rd tab) into a new "2022-Q1" tab placed right
#    after it (and therefore right before "总计").
# ---------------------------------------------------------------------
$q4 = $wb.Worksheets.Item(3)
$q4.Copy([System.Reflection.Missing]::Value, $q4)
$q1new = $wb.Worksheets.Item(4)
$q1new.Name = "2022-Q1"

# ---------------------------------------------------------------------
# 2) Overwrite the cloned data with the 2022-Q1 numbers.
#    Header row (row 1) and column D label already read "基金规模" on
#    the cloned sheet, so only the data rows (2-4) need to change.
#    Numeric-looking text fields are entered with a leading quote so
#    they stay text (matching the source data, which stores these as
#    strings, not numbers).
# ---------------------------------------------------------------------

# Row 2 - 519678 银河消费驱动混合
$q1new.Range("A2").Value = 0
$q1new.Range("B2").Value = "'519678"
$q1new.Range("C2").Value = "银河消费驱动混合"
$q1new.Range("D2").Value = "'1.06"
$q1new.Range("E2").Value = "'75.49"
$q1new.Range("F2").Value = "'8.38"
$q1new.Range("G2").Value = "'0.0888"
$q1new.Range("H2").Value = 2

# Row 3 - 519629 银河睿利灵活配置混合A
$q1new.Range("A3").Value = 1
$q1new.Range("B3").Value = "'519629"
$q1new.Range("C3").Value = "银河睿利灵活配置混合A"
$q1new.Range("D3").Value = "'2.24"
$q1new.Range("E3").Value = "'24.70"
$q1new.Range("F3").Value = "'1.50"
$q1new.Range("G3").Value = "'0.0336"
$q1new.Range("H3").Value = 7

# Row 4 - 519630 银河睿利灵活配置混合C (new row, doesn't exist on the
# cloned sheet yet - copy the format down from row 3 first, then fill
# in the values).
$q1new.Range("A3:H3").Copy()
$q1new.Range("A4").PasteSpecial(-4122)

$q1new.Range("A4").Value = 2
$q1new.Range("B4").Value = "'519630"
$q1new.Range("C4").Value = "银河睿利灵活配置混合C"
$q1new.Range("D4").Value = "'1.94"
$q1new.Range("E4").Value = "'24.70"
$q1new.Range("F4").Value = "'1.50"
$q1new.Range("G4").Value = "'0.0291"
$q1new.Range("H4").Value = 7

# ---------------------------------------------------------------------
# 3) "总计" tab: insert a new row 2 for 2022-Q1 and renumber the index
#    column (A) for every row so it stays 0,1,2,3.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item(5)

$total.Rows("2:2").Insert(-4121)
$total.Range("B2:D2").ClearFormats()

# match the style used by the rest of the index column (A)
$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122)

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 3
$total.Range("D2").Value = 0.15

$total.Range("A3").Value = 1
$total.Range("A4").Value = 2
$total.Range("A5").Value = 3

# Restore the originally-active tab (sheet activation above would
# otherwise leave "2022-Q1" selected).
$wb.Worksheets.Item(1).Activate()
